$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1984.25
$ws.Range("I17").Value = 993.75
$ws.Range("K17").Value = 2981.25
$ws.Range("M17").Value = -2813.25
# Row 33
$ws.Range("H33").Value = 541.82355
$ws.Range("I33").Value = 393.16666
$ws.Range("J33").Value = 898.6
$ws.Range("K33").Value = 393.16666
$ws.Range("L33").Value = 898.6
$ws.Range("M33").Value = -164.16666
$ws.Range("N33").Value = -1356.6
# Row 42
$ws.Range("H42").Value = 1000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
# Row 43
$ws.Range("H43").Value = 2560.875
$ws.Range("J43").Value = 3081.1667
$ws.Range("L43").Value = 3081.1667
$ws.Range("N43").Value = -3219.1667
# Row 47
$ws.Range("H47").Value = 122270
$ws.Range("J47").Value = 122270
$ws.Range("L47").Value = 122270
$ws.Range("N47").Value = -124214
# Row 51
$ws.Range("H51").Value = 3700
# Row 108
$ws.Range("H108").Value = 87180
$ws.Range("J108").Value = 87180
$ws.Range("L108").Value = 87180
$ws.Range("N108").Value = -94860
# Row 109
$ws.Range("H109").Value = 77378.57000000001
$ws.Range("J109").Value = 77378.57000000001
$ws.Range("L109").Value = 77378.57000000001
$ws.Range("N109").Value = -80152.57000000001
# Row 110
$ws.Range("H110").Value = 63635.715
$ws.Range("J110").Value = 63635.715
$ws.Range("L110").Value = 63635.715
$ws.Range("N110").Value = -71815.715
# Row 114
$ws.Range("H114").Value = 91463.25
$ws.Range("J114").Value = 91463.25
$ws.Range("L114").Value = 91463.25
$ws.Range("N114").Value = -100141.25
# Row 117
$ws.Range("H117").Value = 90405.57000000001
$ws.Range("J117").Value = 90405.57000000001
$ws.Range("L117").Value = 90405.57000000001
$ws.Range("N117").Value = -99583.57000000001
# Row 120
$ws.Range("H120").Value = 41044.832
$ws.Range("J120").Value = 41044.832
$ws.Range("L120").Value = 41044.832
$ws.Range("N120").Value = -50720.832
# Row 132
$ws.Range("H132").Value = 2576.9333
$ws.Range("I132").Value = 2243
$ws.Range("J132").Value = 3495.25
$ws.Range("K132").Value = 6729
$ws.Range("L132").Value = 10485.75
$ws.Range("M132").Value = -4199
$ws.Range("N132").Value = -15545.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 570785
$ws.Range("J6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("N6").Value = -2346
# Row 140
$ws.Range("H140").Value = 104049
$ws.Range("J140").Value = 104049
$ws.Range("L140").Value = 104049
$ws.Range("N140").Value = -114409

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 1632349.8
$ws.Range("I7").Value = 3260159.8
$ws.Range("J7").Value = 4539.8
$ws.Range("K7").Value = 3260159.8
$ws.Range("L7").Value = 4539.8
$ws.Range("M7").Value = -3260046.8
$ws.Range("N7").Value = -4765.8
# Row 124
$ws.Range("H124").Value = 48783.332
$ws.Range("J124").Value = 48783.332
$ws.Range("L124").Value = 48783.332
$ws.Range("N124").Value = -58603.332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 2417848
$ws.Range("I4").Value = 4080
$ws.Range("K4").Value = 4080
$ws.Range("M4").Value = -3968
# Row 31
$ws.Range("H31").Value = 4162.5454
$ws.Range("I31").Value = 3199
$ws.Range("K31").Value = 3199
$ws.Range("M31").Value = -2904
# Row 34
$ws.Range("H34").Value = 4162.5454
$ws.Range("I34").Value = 3199
$ws.Range("K34").Value = 3199
$ws.Range("M34").Value = -2997
# Row 132
$ws.Range("H132").Value = 1388.5
$ws.Range("I132").Value = 1376.1111
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 4128.3333
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -1598.3333
$ws.Range("N132").Value = -9560
# Row 140
$ws.Range("H140").Value = 69500
$ws.Range("J140").Value = 69500
$ws.Range("L140").Value = 69500
$ws.Range("N140").Value = -79860

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 999999.5
$ws.Range("I8").Value = 999999.5
$ws.Range("K8").Value = 2999998.5
$ws.Range("M8").Value = -2999859.5
# Row 33
$ws.Range("H33").Value = 397.8125
$ws.Range("I33").Value = 275
$ws.Range("J33").Value = 415.35715
$ws.Range("K33").Value = 1650
$ws.Range("L33").Value = 2492.1429
$ws.Range("M33").Value = -1367
$ws.Range("N33").Value = -3058.1429
# Row 50
$ws.Range("H50").Value = 359.6
$ws.Range("I50").Value = 299
$ws.Range("K50").Value = 897
$ws.Range("M50").Value = -416
# Row 53
$ws.Range("H53").Value = 359.6
$ws.Range("I53").Value = 299
$ws.Range("K53").Value = 897
$ws.Range("M53").Value = -416
# Row 122
$ws.Range("H122").Value = 2020896.8
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050
# Row 129
$ws.Range("H129").Value = 3321.3333
$ws.Range("J129").Value = 3679.8
$ws.Range("L129").Value = 11039.4
$ws.Range("N129").Value = -21039.4
# Row 132
$ws.Range("H132").Value = 2662.7693
$ws.Range("J132").Value = 3776.6667
$ws.Range("L132").Value = 33990.0003
$ws.Range("N132").Value = -39050.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 7000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 7000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 7000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -7224
# Row 33
$ws.Range("H33").Value = 2010757
$ws.Range("J33").Value = 13446.25
$ws.Range("L33").Value = 13446.25
$ws.Range("N33").Value = -13950.25
# Row 53
$ws.Range("H53").Value = 17500
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 17500
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 17500
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -18762
# Row 127
$ws.Range("H127").Value = 69000
$ws.Range("J127").Value = 69000
$ws.Range("L127").Value = 69000
$ws.Range("N127").Value = -78920

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Range("H74").Value = 58752
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 71669.336
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 71669.336
$ws.Range("M74").Value = -19002
$ws.Range("N74").Value = -73665.336
# Row 77
$ws.Range("H77").Value = 58752
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 71669.336
$ws.Range("K77").Value = 60000
$ws.Range("L77").Value = 215008.008
$ws.Range("M77").Value = -55008
$ws.Range("N77").Value = -224992.008
# Row 97
$ws.Range("H97").Value = 45000
$ws.Range("J97").Value = 45000
$ws.Range("L97").Value = 45000
$ws.Range("N97").Value = -46982

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 35000
$ws.Range("I2").Value = 35000
$ws.Range("K2").Value = 35000
$ws.Range("M2").Value = -34888
# Row 4
$ws.Range("H4").Value = 2680
$ws.Range("J4").Value = 3325
$ws.Range("L4").Value = 3325
$ws.Range("N4").Value = -3551
# Row 94
$ws.Range("H94").Value = 21977.4
$ws.Range("I94").Value = 15890
$ws.Range("J94").Value = 23499.25
$ws.Range("K94").Value = 15890
$ws.Range("L94").Value = 23499.25
$ws.Range("M94").Value = -14989
$ws.Range("N94").Value = -25301.25
# Row 126
$ws.Range("H126").Value = 2867.5
$ws.Range("I126").Value = 2485
$ws.Range("K126").Value = 7455
$ws.Range("M126").Value = -4985
